# This script reworks the "8800014" course sheet so that the content that
# previously lived in rows 13-21 is shifted down one row (to 14-22), with
# two brand-new rows of data inserted (the Portuguese "Objetivos" text in
# row 10, "Docentes responsaveis" value in row 13, and a brand-new
# "Bibliografia" content row 22), matching the course's updated syllabus.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 - Objetivos: fill in the (until now missing) Portuguese objectives
# text in columns B/C (previously these erroneously held the teacher's name).
$ws.Range("B10").Value = 'Levar os alunos a entender as principais questões da gestão da inovação, suas características e pontos críticos para obter o sucesso, por meio de problemas reais de empresas. As atividades serão realizadas em equipe e serão focadas no desenvolvimento das competências necessárias para gerenciar com sucesso todo o processo de inovação, de sua concepção até colocar o produto no mercado.'
$ws.Range("C10").Value = 'Levar os alunos a entender as principais questões da gestão da inovação, suas características e pontos críticos para obter o sucesso, por meio de problemas reais de empresas. As atividades serão realizadas em equipe e serão focadas no desenvolvimento das competências necessárias para gerenciar com sucesso todo o processo de inovação, de sua concepção até colocar o produto no mercado.'

# Row 13 - was "Programa resumido: / Semestral"; becomes the
# "Docentes responsaveis" value row (label stays on row 12), with no label
# of its own and the default (non-custom) row height.
$ws.Range("A13").Clear()
$ws.Range("B13").Value = '849935 - Humberto Felipe da Silva'
$ws.Range("C13").Value = '849935 - Humberto Felipe da Silva'
$ws.Rows.Item(13).AutoFit()

# Row 14 - becomes "Programa resumido:" with its real summary text.
$ws.Range("A14").Value = 'Programa resumido:'
$ws.Range("B14").Value = 'Gestão da inovação. Custos da inovação. Processo de implementação da inovação. Transformando a ideação em negócio. O Mercado de inovação.'
$ws.Range("C14").Value = 'Gestão da inovação. Custos da inovação. Processo de implementação da inovação. Transformando a ideação em negócio. O Mercado de inovação.'

# Row 15 - becomes "Short syllabus:" (shifted down from row 14), now at
# 60pt custom row height.
$ws.Range("A15").Value = 'Short syllabus:'
$ws.Range("B15").Value = 'Innovation management. Costs of innovation. Innovation implementation process. Transforming ideation into business. The Innovation Market'
$ws.Range("C15").Value = 'Innovation management. Costs of innovation. Innovation implementation process. Transforming ideation into business. The Innovation Market'
$ws.Rows.Item(15).RowHeight = 60

# Row 16 - becomes "Programa:" with its full Portuguese syllabus text
# (previously this row held the English "Syllabus:" text).
$ws.Range("A16").Value = 'Programa:'
$ws.Range("B16").Value = 'ProgramaO que leva algumas organizações a terem necessidade de gerenciar a inovação; quais são os principais fatores que impulsionam a inovação e como ocorre o processo de difusão. Quais são os principais tipos de inovação que precisam ser considerados; Quais são as principais estratégias para implantar uma de inovação no mercado; Como e de que maneira as empresas procuram obter ideias inovadoras; Quais e como diferentes fatores influenciam a maneira como os gerentes priorizam as escolhas de inovação; Quais os principais desafios para colocar as inovações em prática; De que forma os gestores podem construir uma organização focada na inovação como estratégia de mercado; quais os principais sistemas de avaliação de sucessos de uma ideia inovadora; Gestão de recursos e programas de inovação em uma empresa.'
$ws.Range("C16").Value = 'ProgramaO que leva algumas organizações a terem necessidade de gerenciar a inovação; quais são os principais fatores que impulsionam a inovação e como ocorre o processo de difusão. Quais são os principais tipos de inovação que precisam ser considerados; Quais são as principais estratégias para implantar uma de inovação no mercado; Como e de que maneira as empresas procuram obter ideias inovadoras; Quais e como diferentes fatores influenciam a maneira como os gerentes priorizam as escolhas de inovação; Quais os principais desafios para colocar as inovações em prática; De que forma os gestores podem construir uma organização focada na inovação como estratégia de mercado; quais os principais sistemas de avaliação de sucessos de uma ideia inovadora; Gestão de recursos e programas de inovação em uma empresa.'

# Row 17 - gains the (previously missing) "Syllabus:" label plus its
# English text in columns B/C, and a 120pt custom row height.
$ws.Range("A17").Value = 'Syllabus:'
$ws.Range("B17").Value = 'Which leads some organizations to have the necessity to manage the innovation; what are the main factors driving innovation and how does the diffusion process occur. What are the main types of innovation that need to be considered; What are the main strategies to implement one of innovation in the market; what and how companies seek innovative ideas; What and how different factors influence the way managers to prioritize innovation choices; What are the key challenges in putting innovations into practice? How managers can build an organization focused on innovation as a market strategy; what are the main success evaluation systems of an innovative idea; Resource management and innovative programs in a company.'
$ws.Range("C17").Value = 'Which leads some organizations to have the necessity to manage the innovation; what are the main factors driving innovation and how does the diffusion process occur. What are the main types of innovation that need to be considered; What are the main strategies to implement one of innovation in the market; what and how companies seek innovative ideas; What and how different factors influence the way managers to prioritize innovation choices; What are the key challenges in putting innovations into practice? How managers can build an organization focused on innovation as a market strategy; what are the main success evaluation systems of an innovative idea; Resource management and innovative programs in a company.'
$ws.Rows.Item(17).RowHeight = 120

# Row 18 - becomes just the bare "Avaliação:" label (B/C cleared out,
# default row height) - previously held "Método:" plus the teacher name.
$ws.Range("A18").Value = 'Avaliação:'
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
$ws.Rows.Item(18).AutoFit()

# Row 19 - becomes "Método:" (its B/C text, the seminar description,
# was already in place one row up and needs no change here).
$ws.Range("A19").Value = 'Método:'

# Row 20 - becomes "Critério:" (B/C text unchanged in value/position).
$ws.Range("A20").Value = 'Critério:'

# Row 21 - becomes "Norma de recuperação:" and shrinks from a 120pt to a
# 60pt custom row height.
$ws.Range("A21").Value = 'Norma de recuperação:'
$ws.Rows.Item(21).RowHeight = 60

# Row 22 - brand new "Bibliografia:" row with its full reference list,
# at a 120pt custom row height.
$ws.Range("A22").Value = 'Bibliografia:'
$ws.Range("B22").Value = 'Gestão de Negócios: Visões e dimensões empresariais da Organização. Autores: Cruz Jr, J.B., Rocha, J.A.O. e Tachizawa, T.Editora: ATLASGestão Empresarial - de Taylor aos nossos diasAutores: Pereira, M. I.  Autor: Ferreira, A. A. e Reis, A.C. F Editora: THOMSON PIONEIRAGestão da inovação: a economia da tecnologia no BrasilAutor: Tigre, P. B.Editora: ElsevierTextos disponibilizados pelo professor da disciplinaArtigos extraídos de revistas especializadas na área de gestão e inovação.'
$ws.Range("C22").Value = 'Gestão de Negócios: Visões e dimensões empresariais da Organização. Autores: Cruz Jr, J.B., Rocha, J.A.O. e Tachizawa, T.Editora: ATLASGestão Empresarial - de Taylor aos nossos diasAutores: Pereira, M. I.  Autor: Ferreira, A. A. e Reis, A.C. F Editora: THOMSON PIONEIRAGestão da inovação: a economia da tecnologia no BrasilAutor: Tigre, P. B.Editora: ElsevierTextos disponibilizados pelo professor da disciplinaArtigos extraídos de revistas especializadas na área de gestão e inovação.'
$ws.Rows.Item(22).RowHeight = 120

# Column B cells that are created for the first time in a row that already
# had a styled column-A cell inherit column A's style instead of column B's
# default (wrap-text) style. Fix that up by pasting the formatting from a
# known-good column B cell onto each newly created column B cell.
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B22").PasteSpecial(-4122)
$excel.CutCopyMode = 0
